$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18; this shifts the existing rows 18-47
# down to 19-48, preserving their data/formatting (as in the diff).
$ws.Rows.Item(18).Insert()

# Populate the new row 18 with the new weekly price record.
$newDate = Get-Date -Year 2023 -Month 1 -Day 27 -Hour 0 -Minute 0 -Second 0

$ws.Cells.Item(18, 1).Value = 4
$ws.Cells.Item(18, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(18, 3).Value = "Los Lagos"
$ws.Cells.Item(18, 4).Value = $newDate
$ws.Cells.Item(18, 5).Value = 10
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100101
$ws.Cells.Item(18, 8).Value = "Berries"
$ws.Cells.Item(18, 9).Value = 100101001
$ws.Cells.Item(18, 10).Value = "Arándano (blue)"
$ws.Cells.Item(18, 11).Value = "Sin especificar"
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 400
$ws.Cells.Item(18, 14).Value = 2000
$ws.Cells.Item(18, 15).Value = 2200
$ws.Cells.Item(18, 16).Value = 2100
$ws.Cells.Item(18, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(18, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(18, 19).Value = 1050
$ws.Cells.Item(18, 20).Value = 2
